$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the commit diff.
# Cells whose new text looks like a plain number (e.g. "593.46") would be
# auto-converted to a numeric value by Excel when assigned via .Value, so for
# those we temporarily force a Text number format, assign the value, then
# ClearFormats() to drop the temporary formatting (restores default style)
# while keeping the cells data type as Text.

$ws.Range("D2").Value = '64.430.64'
$ws.Range("E2").Value = '  +1.18%  '
$ws.Range("D3").Value = '3.162.12'
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.46'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.49'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.68%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").Value = '3.154.63'
$ws.Range("E8").Value = '  +1.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.532'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -0.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.67%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.91'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +3.57%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.463'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.62'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").Value = '3.680.95'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("E16").Value = '  -1.14%  '
$ws.Range("E17").Value = '  +1.11%  '
$ws.Range("D18").Value = '64.179.88'
$ws.Range("E18").Value = '  +0.85%  '
$ws.Range("D19").Value = '3.156.43'
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.19'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("E21").Value = '  +0.69%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.737'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.53'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.14'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +7.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.50'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.70'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +8.49%  '
$ws.Range("B29").Value = 'ImmutableX'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.28'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.23%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.74'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +1.19%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.44'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +8.13%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.55'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.111'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.39%  '
$ws.Range("D35").Value = '0.0₃0850'
$ws.Range("E35").Value = '  -1.84%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.08'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +2.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.25'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +2.76%  '
$ws.Range("E38").Value = '  -1.38%  '
$ws.Range("E39").Value = '  -2.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '468.58'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +4.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.82'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.25'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +5.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.301'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +7.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0376'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.90%  '
$ws.Range("D45").Value = '2.936.01'
$ws.Range("E45").Value = '  +1.46%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.74'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +12.03%  '
$ws.Range("E47").Value = '  -2.27%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '129.39'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.64%  '
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.26'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.17%  '
$ws.Range("E51").Value = '  -0.41%  '
